$d = $word.ActiveDocument

# The paragraph currently last in the body is the one ending with
# "...trabajar en el juego en si..." -- give it a bottom paragraph border,
# then append three new paragraphs of journal text after it.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)

$borders = $lastPara.Borders
$borders.DistanceFromBottom = 1
$bottom = $borders.Item(-3)
$bottom.LineStyle = 1
$bottom.LineWidth = 3
$bottom.ColorIndex = 0

$r = $lastPara.Range
$r.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p1.Range.InsertBefore("7/8")
$p1.Borders.Item(-3).LineStyle = 0

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p2.Range.InsertBefore("Se acomodo el aspecto general de la pagina, se coloco algunos background en los escenarios del juego. Hay que revisar nueva_partida.html ya que no se muestra el backgroud-image" + [char]0x2026)
$p2.Borders.Item(-3).LineStyle = 0

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($d.Paragraphs.Count)
$p3.Range.InsertBefore("Se genero la parte de cargar_partida el .js, .html, .php, el js contiene las funciones las cuales hay que revisar. El html aparece pero falta mostrar imagen, se tiene que verificar si se debe modificar el html para que llame a una funci" + [char]0x00F3 + "n o modificar el js para que muestre los datos. El php se encarga de conectar la tabla " + [char]0x201C + "personaje" + [char]0x201D + " de la base de datos.")
$p3.Borders.Item(-3).LineStyle = 0

Write-Output "done"
